$d = $word.ActiveDocument

# Update the resolution text: "causes an underlying problem" -> "shows the underlying problem"
$d.Content.Find.Execute("causes an underlying problem", $true, $false, $false, $false, $false,
                         $true, 1, $false, "shows the underlying problem", 2)

# Update the date: " - 11/01" -> " - 11/13"
$d.Content.Find.Execute(" – 11/01", $true, $false, $false, $false, $false,
                         $true, 1, $false, " – 11/13", 2)
